$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.328456044197083
$ws.Range("B1").Value = 3.165414333343506
$ws.Range("C1").Value = 5.379557132720947
$ws.Range("D1").Value = 1.738052368164062
$ws.Range("E1").Value = 0.9919607043266296
